$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$lq = [char]0x201C
$rq = [char]0x201D

# ---------------------------------------------------------------------------
# 1) "Joey " + proofErr-wrapped "Ampfer" -> single run "Joey Ampfer"
# ---------------------------------------------------------------------------
$pJoey = $d.Paragraphs.Item(4)
$xmlJoey = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Joey Ampfer</w:t></w:r></w:p>"
[void]$pJoey.Range.InsertXML($xmlJoey)

# ---------------------------------------------------------------------------
# 2) Insert two new bullet items ("Cookie Settings Tab Blocks Content" /
#    detail / detail) right before the "Navigation Menu Blocks Content" item.
# ---------------------------------------------------------------------------
$pNav = $d.Paragraphs.Item(20)

# -- new title bullet (level 0)
$pNav.Range.InsertParagraphBefore()
$pCookieTitle = $d.Paragraphs.Item(20)
$xmlCookieTitle = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:t>Cookie Settings Tab Blocks Content</w:t></w:r></w:p>"
[void]$pCookieTitle.Range.InsertXML($xmlCookieTitle)

# -- new detail bullet (level 1) describing the problem
$pNav = $d.Paragraphs.Item(21)
$pNav.Range.InsertParagraphBefore()
$pCookieDetail1 = $d.Paragraphs.Item(21)
$text1 = "There was a " + $lq + "Cookie Settings" + $rq + " tab on the home page that blocked content from the user."
$xmlCookieDetail1 = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:t>$text1</w:t></w:r></w:p>"
[void]$pCookieDetail1.Range.InsertXML($xmlCookieDetail1)

# -- new detail bullet (level 1) describing the fix, split across 5 runs
$pNav = $d.Paragraphs.Item(22)
$pNav.Range.InsertParagraphBefore()
$pCookieDetail2 = $d.Paragraphs.Item(22)
$xmlCookieDetail2 = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='5'/></w:numPr></w:pPr>" + `
    "<w:r><w:t>This tab has been replaced with a one</w:t></w:r>" + `
    "<w:r><w:t>-</w:t></w:r>" + `
    "<w:r><w:t>time popup that will ask the user to accept or deny cookies</w:t></w:r>" + `
    "<w:r><w:t xml:space='preserve'> and a link in the footer </w:t></w:r>" + `
    "<w:r><w:t>to bring it back up.</w:t></w:r>" + `
    "</w:p>"
[void]$pCookieDetail2.Range.InsertXML($xmlCookieDetail2)

# ---------------------------------------------------------------------------
# 3) "Mobile Viewing Adds Space on the Side" title picks up the
#    <w:lastRenderedPageBreak/> that used to sit on the "Color Inconsistency"
#    detail paragraph; that paragraph loses it.
# ---------------------------------------------------------------------------
$pMobileTitle = $d.Paragraphs.Item(26)
$xmlMobileTitle = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Mobile Viewing Adds Space on the Side</w:t></w:r></w:p>"
[void]$pMobileTitle.Range.InsertXML($xmlMobileTitle)

$pMobileDetail1 = $d.Paragraphs.Item(27)
$xmlMobileDetail1 = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:t>On the old home page, mobile viewing would align the contents to the left. This created an odd space to the right of the information.</w:t></w:r></w:p>"
[void]$pMobileDetail1.Range.InsertXML($xmlMobileDetail1)

# ---------------------------------------------------------------------------
# 4) "Our website centers content on mobile to look" / " more professional."
#    -- drop the grammar-check proofErr wrapping around "mobile" and merge
#    the runs down to two.
# ---------------------------------------------------------------------------
$pMobileDetail2 = $d.Paragraphs.Item(28)
$xmlMobileDetail2 = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:t>Our website centers content on mobile to look</w:t></w:r><w:r><w:t xml:space='preserve'> more professional.</w:t></w:r></w:p>"
[void]$pMobileDetail2.Range.InsertXML($xmlMobileDetail2)

# ---------------------------------------------------------------------------
# 5) "Color Inconsistency with Logo" title / detail -- detail loses the
#    lastRenderedPageBreak (moved above in step 3).
# ---------------------------------------------------------------------------
$pColorDetail = $d.Paragraphs.Item(30)
$xmlColorDetail = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='5'/></w:numPr></w:pPr><w:r><w:t>The old website used a color scheme that did not fit with the simple grey, red, white, and black logo.</w:t></w:r></w:p>"
[void]$pColorDetail.Range.InsertXML($xmlColorDetail)
